# Update the existing workbook to match the new data snapshot.
# - Row 170 gets revised values (B/C/D/E/F).
# - Rows 171-174 are appended with new data (dates 46028-46031).
# - The used range/dimension grows from A1:F170 to A1:F174.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OS_GERAL")

# --- Row 170: corrected values ---
$ws.Cells.Item(170, 2).Value = 2591
$ws.Cells.Item(170, 3).Value = 877
$ws.Cells.Item(170, 4).Value = 265
$ws.Cells.Item(170, 5).Value = 1449
$ws.Cells.Item(170, 6).Value = 573

# --- New rows 171-174 ---
$newRows = @(
    @(46028, 3321, 1048, 485, 1788, 366),
    @(46029, 3235, 826, 656, 1753, 515),
    @(46030, 3366, 828, 746, 1792, 432),
    @(46031, 1387, 106, 565, 716, 0)
)

$r = 171
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # Match the date-time style used by column A in existing data rows (style index 2 -> numFmtId 165)
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $r++
}
